# The source data table originally listed years 2008-2020 in rows 2-14.
# This update drops the two oldest years (2008年, 2009年) and appends a
# new row for 2021年, so the table now spans years 2010-2021 in rows 2-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 2008年 row (row 2); everything below shifts up one row.
$ws.Rows.Item(2).Delete()
# Remove the (now-shifted) 2009年 row, which is again row 2.
$ws.Rows.Item(2).Delete()

# After the two deletions, the last data row (2020年) sits at row 12.
# Copy its formatting (bold/bordered/centered label style, and the
# blank placeholder in column O) down to the new row 13 before filling
# in the 2021年 figures.
$ws.Range("A12:AH12").Copy()
$ws.Range("A13:AH13").PasteSpecial(-4122)

# Populate the new 2021年 row.
$ws.Cells.Item(13, 1).Value = "2021年"
$ws.Cells.Item(13, 2).Value = 278514971
$ws.Cells.Item(13, 3).Value = 1123624
$ws.Cells.Item(13, 4).Value = 154965764
$ws.Cells.Item(13, 5).Value = 177201520
$ws.Cells.Item(13, 6).Value = 1064411
$ws.Cells.Item(13, 7).Value = 95933
$ws.Cells.Item(13, 8).Value = 1379315
$ws.Cells.Item(13, 9).Value = 10022580
$ws.Cells.Item(13, 10).Value = 529417
$ws.Cells.Item(13, 11).Value = 2008048
$ws.Cells.Item(13, 12).Value = 130216
$ws.Cells.Item(13, 13).Value = 4018796
$ws.Cells.Item(13, 14).Value = 269579
# Column O (15) is intentionally left blank, matching the empty
# inlineStr cells used for this column in every other recent row.
$ws.Cells.Item(13, 16).Value = 22952507
$ws.Cells.Item(13, 17).Value = 1669749
$ws.Cells.Item(13, 18).Value = 69374017
$ws.Cells.Item(13, 19).Value = 6720228
$ws.Cells.Item(13, 20).Value = 312871
$ws.Cells.Item(13, 21).Value = 488287
$ws.Cells.Item(13, 22).Value = 13479438
$ws.Cells.Item(13, 23).Value = 486390
$ws.Cells.Item(13, 24).Value = 991721
$ws.Cells.Item(13, 25).Value = 3443047
$ws.Cells.Item(13, 26).Value = 7717236
$ws.Cells.Item(13, 27).Value = 111552
$ws.Cells.Item(13, 28).Value = 1900531
$ws.Cells.Item(13, 29).Value = 44747734
$ws.Cells.Item(13, 30).Value = 1085211
$ws.Cells.Item(13, 31).Value = 171602
$ws.Cells.Item(13, 32).Value = 4418621
$ws.Cells.Item(13, 33).Value = 660150
$ws.Cells.Item(13, 34).Value = 32287
